$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1: 99.97 -> 0M
$t.Cell(1,1).Range.Text = "0M"

# Row 2: 0.25 -> 0M
$t.Cell(2,1).Range.Text = "0M"

# Row 3: 825 -> 0M
$t.Cell(3,1).Range.Text = "0M"

# Row 4: 626 -> 1391
$t.Cell(4,1).Range.Text = "1391"

# Row 5: 0.00003 -> 0.00002
$t.Cell(5,1).Range.Text = "0.00002"

# Row 6: 0.00254 -> 0.00259
$t.Cell(6,1).Range.Text = "0.00259"

# Row 7: 0.00012 -> 0.00015
$t.Cell(7,1).Range.Text = "0.00015"

# Row 9: 0.00012 -> 0.00023
$t.Cell(9,1).Range.Text = "0.00023"

# Row 10: 0.00014 -> 0.00027
$t.Cell(10,1).Range.Text = "0.00027"

# Row 11: 0.00018 -> 0.00037
$t.Cell(11,1).Range.Text = "0.00037"

# Row 12: 0.07985 -> 0.24562
$t.Cell(12,1).Range.Text = "0.24562"

# Row 44: collapse multi-run tabbed content to single value 99.97
$t.Cell(44,1).Range.Text = "99.97"

# Row 45: collapse multi-run tabbed content to single value 0.25
$t.Cell(45,1).Range.Text = "0.25"

# Row 46: collapse multi-run tabbed content to single value 825
$t.Cell(46,1).Range.Text = "825"
